# Insert a new data row at row 152 (pushing existing rows 152-280 down to
# 153-281) and populate the new row with the additional weekly price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(152).Insert()

$ws.Range("A152").Value = 9
$ws.Range("B152").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C152").Value = "Metropolitana"
$ws.Range("D152").Value = 44729
$ws.Range("E152").Value = 13
$ws.Range("F152").Value = 300000001
$ws.Range("G152").Value = "Rabanito"
$ws.Range("H152").Value = "Sin especificar"
$ws.Range("I152").Value = "Primera"
$ws.Range("J152").Value = 6100
$ws.Range("K152").Value = 2500
$ws.Range("L152").Value = 3000
$ws.Range("M152").Value = 2750
$ws.Range("N152").Value = "$/cien unidades (volumen en unidades)"
$ws.Range("O152").Value = "Provincia de Chacabuco"
$ws.Range("P152").Value = 28
$ws.Range("Q152").Value = 100
$ws.Range("R152").Value = "Hortaliza"
